$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8516770203761439
$ws.Range("C2").Value = 0.2000309853246733
$ws.Range("D2").Value = 0.3391645995911858
$ws.Range("F2").Value = 0.9440443906620857
$ws.Range("G2").Value = 0.3834933228624351
$ws.Range("H2").Value = 0.5511931531123722
$ws.Range("J2").Value = 0.2783003810015288
$ws.Range("M2").Value = 0.4026474531072424
$ws.Range("O2").Value = 1.812125278026031

$ws.Range("B3").Value = 0.7485809239037167
$ws.Range("C3").Value = 0.1771194128600087
$ws.Range("D3").Value = 0.3356954162833716
$ws.Range("F3").Value = 0.9500294777889238
$ws.Range("G3").Value = 0.3870784247122927
$ws.Range("H3").Value = 0.5572397599184526
$ws.Range("J3").Value = 0.27995609549248
$ws.Range("M3").Value = 0.3715158875389051
$ws.Range("O3").Value = 1.832103535816159

$ws.Range("B4").Value = 0.6850901776286662
$ws.Range("C4").Value = 0.1629960467539036
$ws.Range("D4").Value = 0.3336931673584473
$ws.Range("F4").Value = 0.954395047038993
$ws.Range("G4").Value = 0.3896792054913263
$ws.Range("H4").Value = 0.5612835911194978
$ws.Range("J4").Value = 0.2811472355087972
$ws.Range("M4").Value = 0.3524380038455988
$ws.Range("O4").Value = 1.84590088892763

$ws.Range("B5").Value = 0.6591713030798871
$ws.Range("C5").Value = 0.1572271248371919
$ws.Range("D5").Value = 0.3329094898101772
$ws.Range("F5").Value = 0.9563475815632287
$ws.Range("G5").Value = 0.3908392414146107
$ws.Range("H5").Value = 0.5630147324634578
$ws.Range("J5").Value = 0.2816764894352275
$ws.Range("M5").Value = 0.3446734191735956
$ws.Range("O5").Value = 1.851907673487133

$ws.Range("B6").Value = 0.654864773797101
$ws.Range("C6").Value = 0.1562683951849522
$ws.Range("D6").Value = 0.332781312060888
$ws.Range("F6").Value = 0.9566822746984442
$ws.Range("G6").Value = 0.3910379083546118
$ws.Range("H6").Value = 0.5633072141230571
$ws.Range("J6").Value = 0.2817670194566091
$ws.Range("M6").Value = 0.3433847228971203
$ws.Range("O6").Value = 1.852928284607032

$ws.Range("B7").Value = 0.6847408100881864
$ws.Range("C7").Value = 0.162918299276555
$ws.Range("D7").Value = 0.3336824676648575
$ws.Range("F7").Value = 0.9544206771919264
$ws.Range("G7").Value = 0.3896944448009165
$ws.Range("H7").Value = 0.5613066008770602
$ws.Range("J7").Value = 0.2811541956873143
$ws.Range("M7").Value = 0.3523332475668823
$ws.Range("O7").Value = 1.845980343554729

$ws.Range("B8").Value = 0.8161697821508938
$ws.Range("C8").Value = 0.1921428625813633
$ws.Range("D8").Value = 0.3379419540975448
$ws.Range("F8").Value = 0.9459646005735252
$ws.Range("G8").Value = 0.3846463878543531
$ws.Range("H8").Value = 0.5532092648886362
$ws.Range("J8").Value = 0.2788350298495175
$ws.Range("M8").Value = 0.3919058930738188
$ws.Range("O8").Value = 1.818695718370677

$ws.Range("B9").Value = 1.072336168005336
$ws.Range("C9").Value = 0.2489948470000911
$ws.Range("D9").Value = 0.3473051692525218
$ws.Range("F9").Value = 0.9348703126687212
$ws.Range("G9").Value = 0.3779292818056632
$ws.Range("H9").Value = 0.539959779449724
$ws.Range("J9").Value = 0.275673514319287
$ws.Range("M9").Value = 0.4697836251920364
$ws.Range("O9").Value = 1.777364695698125

$ws.Range("B10").Value = 1.259521010608751
$ws.Range("C10").Value = 0.2904673570007503
$ws.Range("D10").Value = 0.354795966851583
$ws.Range("F10").Value = 0.9300766687866684
$ws.Range("G10").Value = 0.3749512870030998
$ws.Range("H10").Value = 0.5318304572574704
$ws.Range("J10").Value = 0.2741983864668285
$ws.Range("M10").Value = 0.5271504761823849
$ws.Range("O10").Value = 1.754461687787654

$ws.Range("B11").Value = 1.344441917236509
$ws.Range("C11").Value = 0.3092663904974984
$ws.Range("D11").Value = 0.3583356954550823
$ws.Range("F11").Value = 0.9286274916458126
$ws.Range("G11").Value = 0.3740251635060474
$ws.Range("H11").Value = 0.5284812154072824
$ws.Range("J11").Value = 0.2737119203217731
$ws.Range("M11").Value = 0.5532772956062786
$ws.Range("O11").Value = 1.745671899058408

$ws.Range("B12").Value = 1.376564665238334
$ws.Range("C12").Value = 0.3163750836798442
$ws.Range("D12").Value = 0.3596950058382475
$ws.Range("F12").Value = 0.9281840992845289
$ws.Range("G12").Value = 0.3737363838104955
$ws.Range("H12").Value = 0.527263147789796
$ws.Range("J12").Value = 0.2735542907708606
$ws.Range("M12").Value = 0.5631747922779908
$ws.Range("O12").Value = 1.742578378205735

$ws.Range("B13").Value = 1.369648040257005
$ws.Range("C13").Value = 0.3148445554956254
$ws.Range("D13").Value = 0.3594014155131617
$ws.Range("F13").Value = 0.9282749019059011
$ws.Range("G13").Value = 0.3737958189799997
$ws.Range("H13").Value = 0.5275232461067105
$ws.Range("J13").Value = 0.2735870560793146
$ws.Range("M13").Value = 0.5610430269818778
$ws.Range("O13").Value = 1.743234160571745

$ws.Range("B14").Value = 1.347085386009837
$ws.Range("C14").Value = 0.309851431830424
$ws.Range("D14").Value = 0.35844714873096
$ws.Range("F14").Value = 0.9285889004202019
$ws.Range("G14").Value = 0.374000162718815
$ws.Range("H14").Value = 0.5283799974934027
$ws.Range("J14").Value = 0.2736984190044893
$ws.Range("M14").Value = 0.5540914947200264
$ws.Range("O14").Value = 1.745412679598815

$ws.Range("B15").Value = 1.333260502620817
$ws.Range("C15").Value = 0.306791671589167
$ws.Range("D15").Value = 0.3578650898742382
$ws.Range("F15").Value = 0.9287949628392198
$ws.Range("G15").Value = 0.3741334020662777
$ws.Range("H15").Value = 0.5289113239421823
$ws.Range("J15").Value = 0.2737700953774294
$ws.Range("M15").Value = 0.5498339646095189
$ws.Range("O15").Value = 1.746777708998991

$ws.Range("B16").Value = 1.253966436916244
$ws.Range("C16").Value = 0.2892374083064055
$ws.Range("D16").Value = 0.3545672867839897
$ws.Range("F16").Value = 0.9301861117321195
$ws.Range("G16").Value = 0.375020460656323
$ws.Range("H16").Value = 0.5320563622364745
$ws.Range("J16").Value = 0.2742338959863631
$ws.Range("M16").Value = 0.5254435926829473
$ws.Range("O16").Value = 1.755068962836461

$ws.Range("B17").Value = 1.205261778958402
$ws.Range("C17").Value = 0.2784509455536011
$ws.Range("D17").Value = 0.3525779556058382
$ws.Range("F17").Value = 0.9312270238202061
$ws.Range("G17").Value = 0.3756746219459757
$ws.Range("H17").Value = 0.5340751248454012
$ws.Range("J17").Value = 0.274565727648806
$ws.Range("M17").Value = 0.510488286956388
$ws.Range("O17").Value = 1.760573113049105

$ws.Range("B18").Value = 1.177226567920798
$ws.Range("C18").Value = 0.2722405762334574
$ws.Range("D18").Value = 0.3514461877933286
$ws.Range("F18").Value = 0.931894567369298
$ws.Range("G18").Value = 0.3760912008868758
$ws.Range("H18").Value = 0.5352690936400322
$ws.Range("J18").Value = 0.2747739591590488
$ws.Range("M18").Value = 0.5018892787616096
$ws.Range("O18").Value = 1.76389223052783

$ws.Range("B19").Value = 1.167730676577776
$ws.Range("C19").Value = 0.2701367870355682
$ws.Range("D19").Value = 0.3510651310187285
$ws.Range("F19").Value = 0.9321324032640064
$ws.Range("G19").Value = 0.3762391638339011
$ws.Range("H19").Value = 0.5356789878041752
$ws.Range("J19").Value = 0.2748474449036209
$ws.Range("M19").Value = 0.4989783144933568
$ws.Range("O19").Value = 1.765042329116753

$ws.Range("B20").Value = 1.210448719639487
$ws.Range("C20").Value = 0.2795998363276055
$ws.Range("D20").Value = 0.3527884362970894
$ws.Range("F20").Value = 0.9311090910157986
$ws.Range("G20").Value = 0.3756008099807175
$ws.Range("H20").Value = 0.5338568259324816
$ws.Range("J20").Value = 0.2745286054883138
$ws.Range("M20").Value = 0.5120800097041354
$ws.Range("O20").Value = 1.759971317958019

$ws.Range("B21").Value = 1.353713547643906
$ws.Range("C21").Value = 0.3113183108098383
$ws.Range("D21").Value = 0.3587269281517393
$ws.Range("F21").Value = 0.9284938099354534
$ws.Range("G21").Value = 0.3739384591017227
$ws.Range("H21").Value = 0.5281269855026096
$ws.Range("J21").Value = 0.2736649871759482
$ws.Range("M21").Value = 0.5561332282988616
$ws.Range("O21").Value = 1.744766412497938

$ws.Range("B22").Value = 1.447140807256631
$ws.Range("C22").Value = 0.3319891306659031
$ws.Range("D22").Value = 0.3627181429775135
$ws.Range("F22").Value = 0.9273988940914464
$ws.Range("G22").Value = 0.3732130872014778
$ws.Range("H22").Value = 0.5246749361097756
$ws.Range("J22").Value = 0.2732555378671648
$ws.Range("M22").Value = 0.5849466769321623
$ws.Range("O22").Value = 1.736199122173474

$ws.Range("B23").Value = 1.397296238576416
$ws.Range("C23").Value = 0.3209622762861102
$ws.Range("D23").Value = 0.3605779208992459
$ws.Range("F23").Value = 0.9279269971501023
$ws.Range("G23").Value = 0.3735670973478022
$ws.Range("H23").Value = 0.5264905561667916
$ws.Range("J23").Value = 0.2734598743770462
$ws.Range("M23").Value = 0.5695665461823438
$ws.Range("O23").Value = 1.740646046197043

$ws.Range("B24").Value = 1.208103810319642
$ws.Range("C24").Value = 0.279080451083928
$ws.Range("D24").Value = 0.3526932408236405
$ws.Range("F24").Value = 0.9311621931802279
$ws.Range("G24").Value = 0.3756340542599048
$ws.Range("H24").Value = 0.5339554149747414
$ws.Range("J24").Value = 0.2745453340434594
$ws.Range("M24").Value = 0.5113603949379097
$ws.Range("O24").Value = 1.760242907667575

$ws.Range("B25").Value = 1.003210874228671
$ws.Range("C25").Value = 0.2336657000380171
$ws.Range("D25").Value = 0.3446644503275564
$ws.Range("F25").Value = 0.9372828651533922
$ws.Range("G25").Value = 0.3794039886594049
$ws.Range("H25").Value = 0.5432624022774917
$ws.Range("J25").Value = 0.2763801121061533
$ws.Range("M25").Value = 0.4486879843561553
$ws.Range("O25").Value = 1.787238234351477
